$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Fgf17/Fgfr2 -> ECs)
$ws.Range("G2").Value = 0.07496433333333334
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.036942
$ws.Range("N2").Value = 0.110826
$ws.Range("O2").Value = 0.02099032928903418
$ws.Range("P2").Value = 0.02099032928903418
$ws.Range("Q2").Value = 0.002769332402
$ws.Range("R2").Value = 0.024923991618
$ws.Range("S2").Value = 0.02099032928903418
$ws.Range("T2").Value = 0.02099032928903418

# Row 3 (FAPs -> Fgf17/Fgfr2 -> FAPs)
$ws.Range("G3").Value = 0.07496433333333334
$ws.Range("O3").Value = 0.5358731102718634
$ws.Range("P3").Value = 0.5358731102718634
$ws.Range("Q3").Value = 0.07069973735055557
$ws.Range("R3").Value = 0.636297636155
$ws.Range("S3").Value = 0.5358731102718634
$ws.Range("T3").Value = 0.5358731102718634

# Row 4 (FAPs -> Fgf17/Fgfr2 -> MuSCs)
$ws.Range("G4").Value = 0.07496433333333334
$ws.Range("O4").Value = 0.4431365604391025
$ws.Range("P4").Value = 0.4431365604391026
$ws.Range("Q4").Value = 0.05846465857855556
$ws.Range("S4").Value = 0.4431365604391025
$ws.Range("T4").Value = 0.4431365604391026
